# Edit Slide 5 ("CONTRIBUTION") of the RESQURE presentation:
#   - "-CSS Styling"      -> "-CSS Styling & adding colorful UI to the pages."
#   - "-Chat Application" -> "-Chat Application & code combining."
#
# Both runs keep their existing character formatting (Alatsi font, sz=3200,
# dk1 fill, ...) because we only replace the characters of the existing run
# in place via TextRange.Characters(start, length), instead of rewriting the
# whole TextFrame (which would otherwise split the edited text into a new
# run with its own rPr).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(5)

for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shp = $s.Shapes.Item($i)
    if (-not $shp.HasTextFrame) {
        continue
    }

    $tr = $shp.TextFrame.TextRange
    $full = $tr.Text

    if ($full -like "*-CSS Styling*") {
        $oldText = "-CSS Styling"
        $idx = $full.IndexOf($oldText)
        $startPos = $idx + 1
        $target = $tr.Characters($startPos, $oldText.Length)
        $target.Text = "-CSS Styling & adding colorful UI to the pages."
    }

    if ($full -like "*-Chat Application*") {
        $oldText = "-Chat Application"
        $idx = $full.IndexOf($oldText)
        $startPos = $idx + 1
        $target = $tr.Characters($startPos, $oldText.Length)
        $target.Text = "-Chat Application & code combining."
    }
}
